$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: Wins / Losses / Ties in AD1:AF1, styled like the other headers (AC1)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Data rows 2-41: Wins=90, Losses=72, Ties=0
for ($r = 2; $r -le 41; $r++) {
    $ws.Cells.Item($r, 30).Value = 90
    $ws.Cells.Item($r, 31).Value = 72
    $ws.Cells.Item($r, 32).Value = 0
}
